$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Add paragraph border (space-only, no line) around the first paragraph
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5

# Update the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# Replace the ID placeholder text and drop the trailing space run in one go
$d.Content.Find.Execute("**ID__AFFARS_5305_topic_7__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5305_204__ID**", 2)
